# Update the organisation website URL in cell B10 and nudge the UI/view
# state to reflect where the author was working when they made the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("B10")

# Update the text itself: "www.stat.kg" -> "www.stat.gov.kg"
$cell.Value = "www.stat.gov.kg"

# Touching the font's theme colour (to its existing value) forces Excel to
# stamp the cell with its own dedicated font/style entry instead of sharing
# the old one -- mirroring what happened in the real edit (the retyped cell
# picked up its own font record) without altering how the cell actually
# looks (still plain, non-bold, theme colour 1, vertical-top, no wrap).
$cell.Font.ThemeColor = 1

# Move the selection to the cell that was edited; this also drops the old
# "topLeftCell" scroll anchor that pointed elsewhere.
$cell.Select()
